$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 29.7437801361084
$ws.Cells.Item(2, 4).Value = 0.3737801361083939
$ws.Cells.Item(2, 5).Value = 0.1397115901492095
$ws.Cells.Item(3, 2).Value = 29.53999999999999
$ws.Cells.Item(3, 3).Value = 29.31538963317871
$ws.Cells.Item(3, 4).Value = -0.2246103668212811
$ws.Cells.Item(3, 5).Value = 0.05044981688359045
$ws.Cells.Item(4, 3).Value = 29.4737377166748
$ws.Cells.Item(4, 4).Value = -0.07626228332519247
$ws.Cells.Item(4, 5).Value = 0.005815935857971929
$ws.Cells.Item(5, 3).Value = 29.6823787689209
$ws.Cells.Item(5, 4).Value = -0.06762123107910156
$ws.Cells.Item(5, 5).Value = 0.004572630892653251
$ws.Cells.Item(6, 3).Value = 29.89393424987793
$ws.Cells.Item(6, 4).Value = 0.05393424987792628
$ws.Cells.Item(6, 5).Value = 0.00290890330989459
$ws.Cells.Item(7, 3).Value = 29.88800621032715
$ws.Cells.Item(7, 4).Value = 0.07800621032714616
$ws.Cells.Item(7, 5).Value = 0.006084968849602965
$ws.Cells.Item(8, 3).Value = 29.92011070251465
$ws.Cells.Item(8, 4).Value = 0.0001107025146467322
$ws.Cells.Item(8, 5).Value = 0.00000001225504674910996
$ws.Cells.Item(9, 3).Value = 30.01376152038574
$ws.Cells.Item(9, 4).Value = 0.03376152038573821
$ws.Cells.Item(9, 5).Value = 0.001139840258756617
$ws.Cells.Item(10, 2).Value = 30.03999999999999
$ws.Cells.Item(10, 3).Value = 30.16930961608887
$ws.Cells.Item(10, 4).Value = 0.1293096160888751
$ws.Cells.Item(10, 5).Value = 0.01672097681305228
$ws.Cells.Item(11, 2).Value = 30.21000000000001
$ws.Cells.Item(11, 3).Value = 30.22289848327637
$ws.Cells.Item(11, 4).Value = 0.01289848327635923
$ws.Cells.Item(11, 5).Value = 0.0001663708708305187
$ws.Cells.Item(12, 3).Value = 30.32418632507324
$ws.Cells.Item(12, 4).Value = 0.1041863250732433
$ws.Cells.Item(12, 5).Value = 0.01085479033226753
$ws.Cells.Item(13, 3).Value = 30.36539268493652
$ws.Cells.Item(13, 4).Value = -0.01460731506347202
$ws.Cells.Item(13, 5).Value = 0.0002133736533635364
$ws.Cells.Item(14, 3).Value = 30.54109001159668
$ws.Cells.Item(14, 4).Value = 0.101090011596682
$ws.Cells.Item(14, 5).Value = 0.01021919044461729
$ws.Cells.Item(15, 3).Value = 30.4333324432373
$ws.Cells.Item(15, 4).Value = -0.04666755676269929
$ws.Cells.Item(15, 5).Value = 0.00217786085419976
$ws.Cells.Item(16, 3).Value = 30.44302177429199
$ws.Cells.Item(16, 4).Value = -0.2469782257080055
$ws.Cells.Item(16, 5).Value = 0.06099824397387453
$ws.Cells.Item(17, 3).Value = 30.52885246276855
$ws.Cells.Item(17, 4).Value = -0.2211475372314453
$ws.Cells.Item(17, 5).Value = 0.04890623322353349
$ws.Cells.Item(18, 3).Value = 30.65643119812012
$ws.Cells.Item(18, 4).Value = -0.2835688018798805
$ws.Cells.Item(18, 5).Value = 0.08041126539959094
$ws.Cells.Item(19, 3).Value = 30.76356315612793
$ws.Cells.Item(19, 4).Value = -0.1864368438720732
$ws.Cells.Item(19, 5).Value = 0.03475869675297978
$ws.Cells.Item(20, 3).Value = 31.10555839538574
$ws.Cells.Item(20, 4).Value = 0.08555839538574617
$ws.Cells.Item(20, 5).Value = 0.007320239020983671
$ws.Cells.Item(21, 3).Value = 31.2354564666748
$ws.Cells.Item(21, 4).Value = 0.1154564666748001
$ws.Cells.Item(21, 5).Value = 0.01333019569702924
$ws.Cells.Item(22, 3).Value = 31.31316184997559
$ws.Cells.Item(22, 4).Value = 0.0331618499755848
$ws.Cells.Item(22, 5).Value = 0.001099708293803194
$ws.Cells.Item(23, 3).Value = 31.2291316986084
$ws.Cells.Item(23, 4).Value = -0.150868301391597
$ws.Cells.Item(23, 5).Value = 0.02276124436478575
$ws.Cells.Item(24, 3).Value = 31.40299797058105
$ws.Cells.Item(24, 4).Value = -0.1770020294189436
$ws.Cells.Item(24, 5).Value = 0.03132971841842458
$ws.Cells.Item(25, 2).Value = 31.65000000000001
$ws.Cells.Item(25, 3).Value = 31.89637565612793
$ws.Cells.Item(25, 4).Value = 0.246375656127924
$ws.Cells.Item(25, 5).Value = 0.06070096393246505
$ws.Cells.Item(26, 3).Value = 32.42288589477539
$ws.Cells.Item(26, 4).Value = 0.5428858947753952
$ws.Cells.Item(26, 5).Value = 0.2947250947460814
$ws.Cells.Item(27, 3).Value = 32.39987564086914
$ws.Cells.Item(27, 4).Value = 0.1198756408691395
$ws.Cells.Item(27, 5).Value = 0.0143701692737869
$ws.Cells.Item(28, 3).Value = 32.49711227416992
$ws.Cells.Item(28, 4).Value = 0.04711227416991903
$ws.Cells.Item(28, 5).Value = 0.00221956637746162
$ws.Cells.Item(29, 2).Value = 32.84999999999999
$ws.Cells.Item(29, 3).Value = 32.75205612182617
$ws.Cells.Item(29, 4).Value = -0.09794387817382244
$ws.Cells.Item(29, 5).Value = 0.009593003271728571
$ws.Cells.Item(30, 2).Value = 32.90000000000001
$ws.Cells.Item(30, 3).Value = 32.95458602905273
$ws.Cells.Item(30, 4).Value = 0.05458602905272869
$ws.Cells.Item(30, 5).Value = 0.002979634567745341
$ws.Cells.Item(31, 2).Value = 33.09999999999999
$ws.Cells.Item(31, 3).Value = 32.89218521118164
$ws.Cells.Item(31, 4).Value = -0.2078147888183537
$ws.Cells.Item(31, 5).Value = 0.04318698645161694
$ws.Cells.Item(32, 2).Value = 33.40000000000001
$ws.Cells.Item(32, 3).Value = 33.6392707824707
$ws.Cells.Item(32, 4).Value = 0.2392707824706974
$ws.Cells.Item(32, 5).Value = 0.05725050734413981
$ws.Cells.Item(33, 3).Value = 33.65039443969727
$ws.Cells.Item(33, 4).Value = -0.04960556030273722
$ws.Cells.Item(33, 5).Value = 0.002460711612948498
$ws.Cells.Item(34, 2).Value = 34.09999999999999
$ws.Cells.Item(34, 3).Value = 33.84643173217773
$ws.Cells.Item(34, 4).Value = -0.2535682678222599
$ws.Cells.Item(34, 5).Value = 0.06429686644638134
$ws.Cells.Item(35, 2).Value = 34.40000000000001
$ws.Cells.Item(35, 3).Value = 34.40230941772461
$ws.Cells.Item(35, 4).Value = 0.002309417724603691
$ws.Cells.Item(35, 5).Value = 0.000005333410226713688
$ws.Cells.Item(36, 2).Value = 34.90000000000001
$ws.Cells.Item(36, 3).Value = 35.04678726196289
$ws.Cells.Item(36, 4).Value = 0.1467872619628849
$ws.Cells.Item(36, 5).Value = 0.02154650027456061
$ws.Cells.Item(37, 3).Value = 35.71496200561523
$ws.Cells.Item(37, 4).Value = 0.4149620056152372
$ws.Cells.Item(37, 5).Value = 0.1721934661042202
$ws.Cells.Item(38, 3).Value = 36.00539016723633
$ws.Cells.Item(38, 4).Value = 0.3053901672363253
$ws.Cells.Item(38, 5).Value = 0.09326315424463072
$ws.Cells.Item(39, 3).Value = 36.00448989868164
$ws.Cells.Item(39, 4).Value = -0.2955101013183565
$ws.Cells.Item(39, 5).Value = 0.08732621998118534
$ws.Cells.Item(40, 3).Value = 36.59141159057617
$ws.Cells.Item(40, 4).Value = -0.2085884094238253
$ws.Cells.Item(40, 5).Value = 0.04350912454596136
$ws.Cells.Item(41, 3).Value = 37.24056625366211
$ws.Cells.Item(41, 4).Value = -0.05943374633788778
$ws.Cells.Item(41, 5).Value = 0.003532370203756389
$ws.Cells.Item(42, 2).Value = 37.90000000000001
$ws.Cells.Item(42, 3).Value = 37.99029922485352
$ws.Cells.Item(42, 4).Value = 0.09029922485350994
$ws.Cells.Item(42, 5).Value = 0.008153950009144748
$ws.Cells.Item(43, 3).Value = 38.42107772827148
$ws.Cells.Item(43, 4).Value = -0.07892227172851562
$ws.Cells.Item(43, 5).Value = 0.006228724974789657
$ws.Cells.Item(44, 2).Value = 38.90000000000001
$ws.Cells.Item(44, 3).Value = 39.00311660766602
$ws.Cells.Item(44, 4).Value = 0.1031166076660099
$ws.Cells.Item(44, 5).Value = 0.01063303477654582
$ws.Cells.Item(45, 2).Value = 39.40000000000001
$ws.Cells.Item(45, 3).Value = 39.53226852416992
$ws.Cells.Item(45, 4).Value = 0.1322685241699162
$ws.Cells.Item(45, 5).Value = 0.0174949624860877
$ws.Cells.Item(46, 2).Value = 39.90000000000001
$ws.Cells.Item(46, 3).Value = 39.5301399230957
$ws.Cells.Item(46, 4).Value = -0.3698600769043026
$ws.Cells.Item(46, 5).Value = 0.1367964764876566
$ws.Cells.Item(47, 2).Value = 40.09999999999999
$ws.Cells.Item(47, 3).Value = 39.9373664855957
$ws.Cells.Item(47, 4).Value = -0.1626335144042912
$ws.Cells.Item(47, 5).Value = 0.02644966000749079
$ws.Cells.Item(48, 2).Value = 40.59999999999999
$ws.Cells.Item(48, 3).Value = 40.45168685913086
$ws.Cells.Item(48, 4).Value = -0.1483131408691349
$ws.Cells.Item(48, 5).Value = 0.02199678775446786
$ws.Cells.Item(49, 2).Value = 40.90000000000001
$ws.Cells.Item(49, 3).Value = 40.73793411254883
$ws.Cells.Item(49, 4).Value = -0.1620658874511776
$ws.Cells.Item(49, 5).Value = 0.02626535187533775
$ws.Cells.Item(50, 2).Value = 41.20000000000001
$ws.Cells.Item(50, 3).Value = 41.31341552734375
$ws.Cells.Item(50, 4).Value = 0.1134155273437401
$ws.Cells.Item(50, 5).Value = 0.01286308184265865
$ws.Cells.Item(51, 3).Value = 41.75338745117188
$ws.Cells.Item(51, 4).Value = 0.253387451171875
$ws.Cells.Item(51, 5).Value = 0.06420520041137934
$ws.Cells.Item(52, 3).Value = 0.1432662963866917
$ws.Cells.Item(52, 5).Value = 1.856198709984518
$ws.Cells.Item(53, 5).Value = 0.03712397419969036
